$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old (collapsed) "_GoBack" bookmark. It currently sits
#    near the end of the "... Es wurde von uns das Spiel Vier Gewinnt
#    ausgewaehlt ..." paragraph. We will recreate it below so that it
#    spans from the very start of the document, through the end of
#    the "Wenn ein Spieler ... Spalte." paragraph (i.e. everything
#    that remains once the final bullet paragraph is removed).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Delete the whole list paragraph containing the citation
#    "Mastering the game of Go with deep neural networks and tree
#    search" (a "Listenabsatz" / numbered list item referencing
#    numId 1). Locate it via Find so the script does not depend on a
#    hard-coded paragraph index.
# ------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Mastering the game of Go with deep neural networks and tree search",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $targetParagraph = $searchRange.Paragraphs(1)
    $targetParagraph.Range.Delete()
}

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark spanning from the document start
#    up to the end of the now-last piece of body text (i.e. through
#    the end of the "Wenn ein Spieler ... Spalte." paragraph).
# ------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $d.Range(0, $d.Content.End))
